$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "keep_col"
$ws.Range("B5").Value = "issue_date,application_time,grade,installment,total_credit_utilized,accounts,state,vector_feature,non_decimal_feature,sentence_feature,y"

$ws.Range("B6").Select()
